$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2610.0833
$ws.Range("I62").Value = 2373.6365
$ws.Range("J62").Value = 2810.1538
$ws.Range("K62").Value = 2373.6365
$ws.Range("L62").Value = 2810.1538
$ws.Range("M62").Value = -1749.6365
$ws.Range("N62").Value = -4058.1538
$ws.Range("H65").Value = 2610.0833
$ws.Range("I65").Value = 2373.6365
$ws.Range("J65").Value = 2810.1538
$ws.Range("K65").Value = 11868.1825
$ws.Range("L65").Value = 14050.769
$ws.Range("M65").Value = -8748.182500000001
$ws.Range("N65").Value = -20290.769
$ws.Range("H98").Value = 960.8182
$ws.Range("I98").Value = 719.75
$ws.Range("J98").Value = 1603.6666
$ws.Range("K98").Value = 719.75
$ws.Range("L98").Value = 1603.6666
$ws.Range("M98").Value = 778.25
$ws.Range("N98").Value = -4599.6666
$ws.Range("H107").Value = 740.3871
$ws.Range("I107").Value = 710.7083
$ws.Range("K107").Value = 710.7083
$ws.Range("M107").Value = 1209.2917
$ws.Range("H122").Value = 960.8182
$ws.Range("I122").Value = 719.75
$ws.Range("J122").Value = 1603.6666
$ws.Range("K122").Value = 2159.25
$ws.Range("L122").Value = 4810.9998
$ws.Range("M122").Value = 290.75
$ws.Range("N122").Value = -9710.9998
$ws.Range("H129").Value = 759.8615
$ws.Range("J129").Value = 798.45
$ws.Range("L129").Value = 2395.35
$ws.Range("N129").Value = -12395.35
$ws.Range("H137").Value = 83071.875
$ws.Range("I137").Value = 101071.73
$ws.Range("J137").Value = 3072.5557
$ws.Range("K137").Value = 303215.19
$ws.Range("L137").Value = 9217.667099999999
$ws.Range("M137").Value = -300665.19
$ws.Range("N137").Value = -14317.6671
$ws.Range("H141").Value = 2974.5386
$ws.Range("I141").Value = 2706.5557
$ws.Range("J141").Value = 3577.5
$ws.Range("K141").Value = 8119.6671
$ws.Range("L141").Value = 10732.5
$ws.Range("M141").Value = -2939.6671
$ws.Range("N141").Value = -21092.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1004.087
$ws.Range("I2").Value = 927.44446
$ws.Range("J2").Value = 1280
$ws.Range("K2").Value = 927.44446
$ws.Range("L2").Value = 1280
$ws.Range("M2").Value = -814.44446
$ws.Range("N2").Value = -1506
$ws.Range("H45").Value = 2405.3076
$ws.Range("I45").Value = 2207.647
$ws.Range("J45").Value = 2778.6667
$ws.Range("K45").Value = 2207.647
$ws.Range("L45").Value = 2778.6667
$ws.Range("M45").Value = -1830.647
$ws.Range("N45").Value = -3532.6667
$ws.Range("H61").Value = 1745.2195
$ws.Range("I61").Value = 1472.7742
$ws.Range("K61").Value = 1472.7742
$ws.Range("M61").Value = -1260.7742
$ws.Range("H110").Value = 1121.3
$ws.Range("I110").Value = 1073.1111
$ws.Range("K110").Value = 1073.1111
$ws.Range("M110").Value = 971.8888999999999
$ws.Range("H116").Value = 1004.087
$ws.Range("I116").Value = 927.44446
$ws.Range("J116").Value = 1280
$ws.Range("K116").Value = 927.44446
$ws.Range("L116").Value = 1280
$ws.Range("M116").Value = 1366.55554
$ws.Range("N116").Value = -5868
$ws.Range("H136").Value = 1745.2195
$ws.Range("I136").Value = 1472.7742
$ws.Range("K136").Value = 4418.3226
$ws.Range("M136").Value = -1868.3226

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1004.087
$ws.Range("I3").Value = 927.44446
$ws.Range("J3").Value = 1280
$ws.Range("K3").Value = 927.44446
$ws.Range("L3").Value = 1280
$ws.Range("M3").Value = -813.44446
$ws.Range("N3").Value = -1508
$ws.Range("H35").Value = 25000
$ws.Range("J35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("N35").Value = -25620
$ws.Range("H99").Value = 1693.9231
$ws.Range("I99").Value = 1864
$ws.Range("J99").Value = 1587.625
$ws.Range("K99").Value = 1864
$ws.Range("L99").Value = 1587.625
$ws.Range("M99").Value = -366
$ws.Range("N99").Value = -4583.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3843.302
$ws.Range("I31").Value = 2053.1904
$ws.Range("J31").Value = 5018.0625
$ws.Range("K31").Value = 2053.1904
$ws.Range("L31").Value = 5018.0625
$ws.Range("M31").Value = -1758.1904
$ws.Range("N31").Value = -5608.0625
$ws.Range("H34").Value = 3843.302
$ws.Range("I34").Value = 2053.1904
$ws.Range("J34").Value = 5018.0625
$ws.Range("K34").Value = 2053.1904
$ws.Range("L34").Value = 5018.0625
$ws.Range("M34").Value = -1851.1904
$ws.Range("N34").Value = -5422.0625
$ws.Range("H41").Value = 20000
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20856
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H94").Value = 3433
$ws.Range("J94").Value = 3895.7
$ws.Range("L94").Value = 3895.7
$ws.Range("N94").Value = -4797.7
$ws.Range("H99").Value = 3856.5833
$ws.Range("I99").Value = 2758.7778
$ws.Range("J99").Value = 7150
$ws.Range("K99").Value = 2758.7778
$ws.Range("L99").Value = 7150
$ws.Range("M99").Value = -1260.7778
$ws.Range("N99").Value = -10146
$ws.Range("H105").Value = 816.6667
$ws.Range("I105").Value = 816.6667
$ws.Range("K105").Value = 816.6667
$ws.Range("M105").Value = 930.3333
$ws.Range("H126").Value = 3856.5833
$ws.Range("I126").Value = 2758.7778
$ws.Range("J126").Value = 7150
$ws.Range("K126").Value = 8276.3334
$ws.Range("L126").Value = 21450
$ws.Range("M126").Value = -5806.3334
$ws.Range("N126").Value = -26390
$ws.Range("H134").Value = 843.5294
$ws.Range("I134").Value = 788.14813
$ws.Range("J134").Value = 1057.1428
$ws.Range("K134").Value = 2364.44439
$ws.Range("L134").Value = 3171.4284
$ws.Range("M134").Value = 170.5556099999999
$ws.Range("N134").Value = -8241.428400000001
$ws.Range("H141").Value = 29600.658
$ws.Range("J141").Value = 29600.658
$ws.Range("L141").Value = 29600.658
$ws.Range("N141").Value = -39960.658

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3493.25
$ws.Range("I64").Value = 983.8
$ws.Range("J64").Value = 5285.7144
$ws.Range("K64").Value = 2951.4
$ws.Range("L64").Value = 15857.1432
$ws.Range("M64").Value = -2681.4
$ws.Range("N64").Value = -16397.1432
$ws.Range("H67").Value = 3493.25
$ws.Range("I67").Value = 983.8
$ws.Range("J67").Value = 5285.7144
$ws.Range("K67").Value = 2951.4
$ws.Range("L67").Value = 15857.1432
$ws.Range("M67").Value = -2015.4
$ws.Range("N67").Value = -17729.1432
$ws.Range("H114").Value = 3400
$ws.Range("I114").Value = 4600
$ws.Range("K114").Value = 13800
$ws.Range("M114").Value = -10546
$ws.Range("H117").Value = 1519.8462
$ws.Range("I117").Value = 1091.6
$ws.Range("J117").Value = 1787.5
$ws.Range("K117").Value = 3274.8
$ws.Range("L117").Value = 5362.5
$ws.Range("M117").Value = 167.2000000000003
$ws.Range("N117").Value = -12246.5
$ws.Range("H131").Value = 716.7
$ws.Range("J131").Value = 729.30927
$ws.Range("L131").Value = 2187.92781
$ws.Range("N131").Value = -12267.92781

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 31160
$ws.Range("J46").Value = 31160
$ws.Range("L46").Value = 31160
$ws.Range("N46").Value = -31472
$ws.Range("H80").Value = 3826.6667
$ws.Range("I80").Value = 3525
$ws.Range("J80").Value = 4171.4287
$ws.Range("K80").Value = 3525
$ws.Range("L80").Value = 4171.4287
$ws.Range("M80").Value = -2527
$ws.Range("N80").Value = -6167.4287
$ws.Range("H83").Value = 3826.6667
$ws.Range("I83").Value = 3525
$ws.Range("J83").Value = 4171.4287
$ws.Range("K83").Value = 17625
$ws.Range("L83").Value = 20857.1435
$ws.Range("M83").Value = -12633
$ws.Range("N83").Value = -30841.1435
$ws.Range("H97").Value = 1241.5807
$ws.Range("I97").Value = 1346.4348
$ws.Range("J97").Value = 940.125
$ws.Range("K97").Value = 1346.4348
$ws.Range("L97").Value = 940.125
$ws.Range("M97").Value = -850.4348
$ws.Range("N97").Value = -1932.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4517.96
$ws.Range("I7").Value = 4413.8335
$ws.Range("K7").Value = 4413.8335
$ws.Range("M7").Value = -4301.8335
$ws.Range("H40").Value = 3888.5557
$ws.Range("I40").Value = 3949.5715
$ws.Range("K40").Value = 3949.5715
$ws.Range("M40").Value = -3813.5715
$ws.Range("H46").Value = 2003.2858
$ws.Range("I46").Value = 2110.5557
$ws.Range("J46").Value = 1810.2
$ws.Range("K46").Value = 2110.5557
$ws.Range("L46").Value = 1810.2
$ws.Range("M46").Value = -1922.5557
$ws.Range("N46").Value = -2186.2
$ws.Range("H68").Value = 2480.4546
$ws.Range("I68").Value = 1899.8
$ws.Range("K68").Value = 1899.8
$ws.Range("M68").Value = -1150.8
$ws.Range("H71").Value = 2480.4546
$ws.Range("I71").Value = 1899.8
$ws.Range("K71").Value = 9499
$ws.Range("M71").Value = -5755
$ws.Range("H122").Value = 1786488.9
$ws.Range("I122").Value = 1964367.8
$ws.Range("K122").Value = 5893103.4
$ws.Range("M122").Value = -5890653.4
$ws.Range("H126").Value = 4517.96
$ws.Range("I126").Value = 4413.8335
$ws.Range("K126").Value = 13241.5005
$ws.Range("M126").Value = -10771.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 42249.75
$ws.Range("J125").Value = 42249.75
$ws.Range("L125").Value = 42249.75
$ws.Range("N125").Value = -52089.75
$ws.Range("H126").Value = 1905.6451
$ws.Range("I126").Value = 1522.1154
$ws.Range("K126").Value = 4566.3462
$ws.Range("M126").Value = -2096.3462
$ws.Range("H132").Value = 831.9268
$ws.Range("I132").Value = 750.3333
$ws.Range("J132").Value = 947.1177
$ws.Range("K132").Value = 2250.9999
$ws.Range("L132").Value = 2841.3531
$ws.Range("M132").Value = 279.0001000000002
$ws.Range("N132").Value = -7901.3531
